$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain "123.45"-style numeric-looking string.
# Force the Text number format first so Excel stores them as text (matching
# the source data, which keeps these as inline/shared strings, not numbers).
$textForceCells = @("D5", "D6", "D20", "D21", "D25", "D27", "D29", "D33", "D44", "D45", "D47")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.619.88'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.594.37'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '211.17'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '0.516'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '1.818.26'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').Value = '1.569.66'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '26.595.78'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '207.39'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '6.85'
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  -4.24%  '
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '145.80'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '7.13'
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').Value = '15.28'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D33').Value = '0.653'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.281.40'
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('D44').Value = '63.50'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '0.920'
$ws.Range('E45').Value = '  +9.95%  '
$ws.Range('D46').Value = '1.730.80'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('D47').Value = '89.69'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('E50').Value = '  +3.43%  '
$ws.Range('E51').Value = '  -1.55%  '

# Restore the default (Normal) style on the text-forced cells so only the
# value changes (no lingering number-format/style difference on these cells).
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
